# Fix word order in the recurring "Dates à utiliser..." sentence:
# "Campagne Constellation du Lion 2022" -> "Campagne 2022 Constellation du Lion"
# This text appears 4 times in the document (all identical), so replace all.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "Dates à utiliser pour la Campagne Constellation du Lion 2022: 14-23 avril, 14-23 mai",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Dates à utiliser pour la Campagne 2022 Constellation du Lion: 14-23 avril, 14-23 mai",
    2
)
